$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 down into rows 3-11 (process bottom-up
# so we don't overwrite source rows before they've been copied).
for ($r = 10; $r -ge 2; $r--) {
    $destRow = $r + 1
    $ws.Range("B$destRow").Value = $ws.Range("B$r").Value()
    $ws.Range("C$destRow").Value = $ws.Range("C$r").Value()
    $ws.Range("D$destRow").Value = $ws.Range("D$r").Value()
    $ws.Range("E$destRow").Value = $ws.Range("E$r").Value()
    $ws.Range("F$destRow").Value = $ws.Range("F$r").Value()
    $ws.Range("G$destRow").Value = $ws.Range("G$r").Value()
}

# Write the newly computed values into row 2.
$ws.Range("B2").Value = 0.2036364321150554
$ws.Range("C2").Value = 0.3089978501498661
$ws.Range("D2").Value = 0.216765395349195
$ws.Range("E2").Value = 0.4655807076643049
$ws.Range("F2").Value = 0.4333807286052367
$ws.Range("G2").Value = 15
